$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Modify Transaction")

# A1: navigation target changed from the JLG group page to the Center page
$ws.Range("A1").Value = "NavigateToCurrentCenterPage"

# B4: re-enter the "06 April 2015" meeting date as an explicit Text-formatted
# value (adds a new numFmtId=49 "@" style and a fresh shared string)
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "06 April 2015"

# B5 keeps displaying "13 April 2015" (its shared string slot shifts down
# automatically once the old "06 April 2015" entry above is freed)

# Selection moved from B12 to C11
[void]$ws.Range("C11").Select()

# Page setup: paper size 9 (A4), portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
